# This script re-aligns the three "file" summary tables (C/D, G/H, K/L) on
# worksheet "master" so that they all start on the same row (row 3) instead
# of being staggered by one row (G/H) and two rows (K/L).
#
# Strategy:
#   1. Read the existing G/H (rows 4-16) and K/L (rows 5-17) values into
#      memory, preserving their typed Value2 representation.
#   2. Fully clear (contents + formats) the old G/H and K/L ranges and
#      remove the merged cells that used to live there.
#   3. Copy the cell formatting (styles) from the now-canonical C/D column
#      (rows 3-15), cell by cell, onto the G/H and K/L ranges, since each
#      table uses an identical per-row style pattern.
#   4. Write the previously captured values back into G/H and K/L, shifted
#      up by 1 row (G/H) and 2 rows (K/L) so that they line up with C/D.
#   5. Re-create the merged header cells at their new locations.
#   6. Fix up the row heights: rows 4, 5, 10 and 11 no longer hold the
#      "tall" header-style cells, so they revert to the default
#      (non-custom) row height, while rows 3 and 9 keep their custom 25pt
#      height untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: capture existing values (typed) before anything is touched.
# ---------------------------------------------------------------------

# G/H table currently occupies rows 4..16
$ghRows = 4..16
$ghCol1 = @{}
$ghCol2 = @{}
foreach ($r in $ghRows) {
    $ghCol1[$r] = $ws.Range("G$r").Value2
    $ghCol2[$r] = $ws.Range("H$r").Value2
}

# K/L table currently occupies rows 5..17
$klRows = 5..17
$klCol1 = @{}
$klCol2 = @{}
foreach ($r in $klRows) {
    $klCol1[$r] = $ws.Range("K$r").Value2
    $klCol2[$r] = $ws.Range("L$r").Value2
}

# ---------------------------------------------------------------------
# Step 2: remove old merges, then fully clear old contents + formats.
# ---------------------------------------------------------------------

$ws.Range("G4:H4").UnMerge() | Out-Null
$ws.Range("G10:H10").UnMerge() | Out-Null
$ws.Range("K5:L5").UnMerge() | Out-Null
$ws.Range("K11:L11").UnMerge() | Out-Null

$ws.Range("G4:H16").Clear() | Out-Null
$ws.Range("K5:L17").Clear() | Out-Null

# ---------------------------------------------------------------------
# Step 3: copy formatting from C/D (rows 3-15) onto the new G/H and K/L
# locations (rows 3-15), cell by cell (avoiding copying merged multi-cell
# ranges in one go, which would otherwise fabricate extra border styles).
# ---------------------------------------------------------------------

foreach ($r in 3..15) {
    $ws.Range("C$r").Copy() | Out-Null
    $ws.Range("G$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $ws.Range("D$r").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $ws.Range("C$r").Copy() | Out-Null
    $ws.Range("K$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $ws.Range("D$r").Copy() | Out-Null
    $ws.Range("L$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 4: write the captured values back, shifted up so everything lines
# up with the C/D table (G/H moves up by 1 row, K/L moves up by 2 rows).
# ---------------------------------------------------------------------

foreach ($r in $ghRows) {
    $destRow = $r - 1
    $v1 = $ghCol1[$r]
    $v2 = $ghCol2[$r]
    if ($null -ne $v1) { $ws.Range("G$destRow").Value = $v1 }
    if ($null -ne $v2) { $ws.Range("H$destRow").Value = $v2 }
}

foreach ($r in $klRows) {
    $destRow = $r - 2
    $v1 = $klCol1[$r]
    $v2 = $klCol2[$r]
    if ($null -ne $v1) { $ws.Range("K$destRow").Value = $v1 }
    if ($null -ne $v2) { $ws.Range("L$destRow").Value = $v2 }
}

# ---------------------------------------------------------------------
# Step 5: re-create merged header cells at the new aligned locations.
# ---------------------------------------------------------------------

$ws.Range("G3:H3").Merge() | Out-Null
$ws.Range("G9:H9").Merge() | Out-Null
$ws.Range("K3:L3").Merge() | Out-Null
$ws.Range("K9:L9").Merge() | Out-Null

# ---------------------------------------------------------------------
# Step 6: fix row heights - rows 4, 5, 10, 11 are no longer "tall" header
# rows, so drop their custom height back to the sheet default.
# ---------------------------------------------------------------------

$ws.Rows(4).AutoFit() | Out-Null
$ws.Rows(5).AutoFit() | Out-Null
$ws.Rows(10).AutoFit() | Out-Null
$ws.Rows(11).AutoFit() | Out-Null

Write-Host "Re-aligned G/H and K/L tables with C/D."
